$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "mu" column (column I), shifting all later columns left by one.
$ws.Columns.Item(9).Delete()

# Update data values (2023-05-13 Eligible for Experiment re-run results).
$ws.Cells.Item(2, 3).Value = 70.95848971701832
$ws.Cells.Item(2, 4).Value = 79561592
$ws.Cells.Item(2, 5).Value = 1126927007.043005
$ws.Cells.Item(2, 6).Value = 0.004220205241383404
$ws.Cells.Item(2, 7).Value = 25.04851854875356
$ws.Cells.Item(2, 8).Value = 20
$ws.Cells.Item(2, 9).Value = $true
$ws.Cells.Item(2, 10).Value = 76923.07692307692
$ws.Cells.Item(2, 11).Value = 109736.707696296
$ws.Cells.Item(2, 12).Value = 18.08138809386727
$ws.Cells.Item(2, 13).Value = 43.12990664262082
$ws.Cells.Item(3, 3).Value = 175.9556230791904
$ws.Cells.Item(3, 4).Value = 79873912
$ws.Cells.Item(3, 5).Value = 232371787.401807
$ws.Cells.Item(3, 6).Value = 0.003491682868916071
$ws.Cells.Item(3, 7).Value = 302.4089996224546
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = $false
$ws.Cells.Item(3, 10).Value = 76923.07692307692
$ws.Cells.Item(3, 11).Value = 97002.55162876946
$ws.Cells.Item(3, 12).Value = 20.45504954955761
$ws.Cells.Item(3, 13).Value = 322.8640491720122
$ws.Cells.Item(4, 3).Value = 118.5748232692446
$ws.Cells.Item(4, 4).Value = 79817080
$ws.Cells.Item(4, 5).Value = 1099128883.62274
$ws.Cells.Item(4, 6).Value = 0.003200154175598993
$ws.Cells.Item(4, 7).Value = 43.0536231732567
$ws.Cells.Item(4, 8).Value = 20
$ws.Cells.Item(4, 9).Value = $true
$ws.Cells.Item(4, 10).Value = 76923.07692307692
$ws.Cells.Item(4, 11).Value = 91467.14202951529
$ws.Cells.Item(4, 12).Value = 21.69294848372683
$ws.Cells.Item(4, 13).Value = 64.74657165698353
$ws.Cells.Item(5, 3).Value = 80.7821316518251
$ws.Cells.Item(5, 4).Value = 79235960
$ws.Cells.Item(5, 5).Value = 833180992.2729816
$ws.Cells.Item(5, 6).Value = 0.0029241020420135
$ws.Cells.Item(5, 7).Value = 38.41212060549256
$ws.Cells.Item(5, 8).Value = 20
$ws.Cells.Item(5, 9).Value = $false
$ws.Cells.Item(5, 10).Value = 76923.07692307692
$ws.Cells.Item(5, 11).Value = 85958.17114104528
$ws.Cells.Item(5, 12).Value = 23.08322726811183
$ws.Cells.Item(5, 13).Value = 61.49534787360439
$ws.Cells.Item(6, 3).Value = 80.74375796642178
$ws.Cells.Item(6, 4).Value = 79761784
$ws.Cells.Item(6, 5).Value = 895042835.8966057
$ws.Cells.Item(6, 6).Value = 0.00430839303353433
$ws.Cells.Item(6, 7).Value = 35.97741875568727
$ws.Cells.Item(6, 8).Value = 20
$ws.Cells.Item(6, 9).Value = $true
$ws.Cells.Item(6, 10).Value = 76923.07692307692
$ws.Cells.Item(6, 11).Value = 111183.5560113591
$ws.Cells.Item(6, 12).Value = 17.84609227462814
$ws.Cells.Item(6, 13).Value = 53.8235110303154
$ws.Cells.Item(7, 3).Value = 111.7396685830363
$ws.Cells.Item(7, 4).Value = 80025976
$ws.Cells.Item(7, 5).Value = 1405892873.14731
$ws.Cells.Item(7, 6).Value = 0.002292505142216274
$ws.Cells.Item(7, 7).Value = 31.80212449706709
$ws.Cells.Item(7, 8).Value = 20
$ws.Cells.Item(7, 9).Value = $true
$ws.Cells.Item(7, 10).Value = 76923.07692307692
$ws.Cells.Item(7, 11).Value = 72219.37576738487
$ws.Cells.Item(7, 12).Value = 27.47451053012403
$ws.Cells.Item(7, 13).Value = 59.27663502719112
$ws.Cells.Item(8, 3).Value = 123.9098299632842
$ws.Cells.Item(8, 4).Value = 79814520
$ws.Cells.Item(8, 5).Value = 1235208275.756846
$ws.Cells.Item(8, 6).Value = 0.0009087328629546882
$ws.Cells.Item(8, 7).Value = 40.03293936701237
$ws.Cells.Item(8, 8).Value = 9
$ws.Cells.Item(8, 9).Value = $false
$ws.Cells.Item(8, 10).Value = 76923.07692307692
$ws.Cells.Item(8, 11).Value = 34408.26517493059
$ws.Cells.Item(8, 12).Value = 57.66614474494507
$ws.Cells.Item(8, 13).Value = 97.69908411195745
$ws.Cells.Item(9, 3).Value = 110.2705246966666
$ws.Cells.Item(9, 4).Value = 79492984
$ws.Cells.Item(9, 5).Value = 967188858.4782858
$ws.Cells.Item(9, 6).Value = 0.002054846491575216
$ws.Cells.Item(9, 7).Value = 45.31551919019817
$ws.Cells.Item(9, 8).Value = 15
$ws.Cells.Item(9, 9).Value = $true
$ws.Cells.Item(9, 10).Value = 76923.07692307692
$ws.Cells.Item(9, 11).Value = 66574.94907215168
$ws.Cells.Item(9, 12).Value = 29.80388310698668
$ws.Cells.Item(9, 13).Value = 75.11940229718485
$ws.Cells.Item(10, 3).Value = 109.2942173926688
$ws.Cells.Item(10, 4).Value = 80046456
$ws.Cells.Item(10, 5).Value = 1828838693.402923
$ws.Cells.Item(10, 6).Value = 0.002065237360215993
$ws.Cells.Item(10, 7).Value = 23.91849755567599
$ws.Cells.Item(10, 8).Value = 20
$ws.Cells.Item(10, 9).Value = $true
$ws.Cells.Item(10, 10).Value = 76923.07692307692
$ws.Cells.Item(10, 11).Value = 66827.82902013171
$ws.Cells.Item(10, 12).Value = 29.69110367781463
$ws.Cells.Item(10, 13).Value = 53.60960123349062
$ws.Cells.Item(11, 3).Value = 153.218164522722
$ws.Cells.Item(11, 4).Value = 79497592
$ws.Cells.Item(11, 5).Value = 788234814.5322708
$ws.Cells.Item(11, 6).Value = 0.002172038545255918
$ws.Cells.Item(11, 7).Value = 77.26425492538016
$ws.Cells.Item(11, 8).Value = 1
$ws.Cells.Item(11, 9).Value = $false
$ws.Cells.Item(11, 10).Value = 76923.07692307692
$ws.Cells.Item(11, 11).Value = 69394.16014641368
$ws.Cells.Item(11, 12).Value = 28.59306886650957
$ws.Cells.Item(11, 13).Value = 105.8573237918897
$ws.Cells.Item(12, 3).Value = 52.69345785456528
$ws.Cells.Item(12, 4).Value = 79762808
$ws.Cells.Item(12, 5).Value = 1304380947.56723
$ws.Cells.Item(12, 6).Value = 0.003557886963972279
$ws.Cells.Item(12, 7).Value = 16.11100717757591
$ws.Cells.Item(12, 8).Value = 20
$ws.Cells.Item(12, 9).Value = $true
$ws.Cells.Item(12, 10).Value = 76923.07692307692
$ws.Cells.Item(12, 11).Value = 98222.0433687905
$ws.Cells.Item(12, 12).Value = 20.20108655803495
$ws.Cells.Item(12, 13).Value = 36.31209373561087
$ws.Cells.Item(13, 3).Value = 185.9418346898805
$ws.Cells.Item(13, 4).Value = 80041848
$ws.Cells.Item(13, 5).Value = 983761988.3588175
$ws.Cells.Item(13, 6).Value = 0.002715096879788233
$ws.Cells.Item(13, 7).Value = 75.64394764793487
$ws.Cells.Item(13, 8).Value = 2
$ws.Cells.Item(13, 9).Value = $false
$ws.Cells.Item(13, 10).Value = 76923.07692307692
$ws.Cells.Item(13, 11).Value = 81597.38144287353
$ws.Cells.Item(13, 12).Value = 24.31685876328196
$ws.Cells.Item(13, 13).Value = 99.96080641121682
$ws.Cells.Item(14, 3).Value = 197.075549015551
$ws.Cells.Item(14, 4).Value = 79997816
$ws.Cells.Item(14, 5).Value = 1645456896.429853
$ws.Cells.Item(14, 6).Value = 0.001254936285397894
$ws.Cells.Item(14, 7).Value = 47.90649193683429
$ws.Cells.Item(14, 8).Value = 9
$ws.Cells.Item(14, 9).Value = $false
$ws.Cells.Item(14, 10).Value = 76923.07692307692
$ws.Cells.Item(14, 11).Value = 45143.10254511117
$ws.Cells.Item(14, 12).Value = 43.95338131705085
$ws.Cells.Item(14, 13).Value = 91.85987325388514
$ws.Cells.Item(15, 3).Value = 72.50100601321395
$ws.Cells.Item(15, 4).Value = 79731064
$ws.Cells.Item(15, 5).Value = 1371947743.31511
$ws.Cells.Item(15, 6).Value = 0.003943618243461979
$ws.Cells.Item(15, 7).Value = 21.06706461186349
$ws.Cells.Item(15, 8).Value = 20
$ws.Cells.Item(15, 9).Value = $true
$ws.Cells.Item(15, 10).Value = 76923.07692307692
$ws.Cells.Item(15, 11).Value = 105072.5329416696
$ws.Cells.Item(15, 12).Value = 18.88402177476308
$ws.Cells.Item(15, 13).Value = 39.95108638662657
$ws.Cells.Item(16, 3).Value = 96.3739219971601
$ws.Cells.Item(16, 4).Value = 79832952
$ws.Cells.Item(16, 5).Value = 789119212.298843
$ws.Cells.Item(16, 6).Value = 0.002360910960038808
$ws.Cells.Item(16, 7).Value = 48.74938139218275
$ws.Cells.Item(16, 8).Value = 15
$ws.Cells.Item(16, 9).Value = $true
$ws.Cells.Item(16, 10).Value = 76923.07692307692
$ws.Cells.Item(16, 11).Value = 73792.20098831299
$ws.Cells.Item(16, 12).Value = 26.88891201814472
$ws.Cells.Item(16, 13).Value = 75.63829341032746
$ws.Cells.Item(17, 3).Value = 132.3194100858378
$ws.Cells.Item(17, 4).Value = 79816568
$ws.Cells.Item(17, 5).Value = 1734880216.672365
$ws.Cells.Item(17, 6).Value = 0.003627136141331431
$ws.Cells.Item(17, 7).Value = 30.4380702809947
$ws.Cells.Item(17, 8).Value = 20
$ws.Cells.Item(17, 9).Value = $true
$ws.Cells.Item(17, 10).Value = 76923.07692307692
$ws.Cells.Item(17, 11).Value = 99483.44712401924
$ws.Cells.Item(17, 12).Value = 19.94494619317365
$ws.Cells.Item(17, 13).Value = 50.38301647416836
$ws.Cells.Item(18, 3).Value = 176.1311528522601
$ws.Cells.Item(18, 4).Value = 79846264
$ws.Cells.Item(18, 5).Value = 1118402169.629412
$ws.Cells.Item(18, 6).Value = 0.004926720339789661
$ws.Cells.Item(18, 7).Value = 62.87279706335826
$ws.Cells.Item(18, 8).Value = 13
$ws.Cells.Item(18, 9).Value = $true
$ws.Cells.Item(18, 10).Value = 76923.07692307692
$ws.Cells.Item(18, 11).Value = 120830.5506456392
$ws.Cells.Item(18, 12).Value = 16.42127747823526
$ws.Cells.Item(18, 13).Value = 79.29407454159352
$ws.Cells.Item(19, 3).Value = 102.2661244817619
$ws.Cells.Item(19, 4).Value = 79854968
$ws.Cells.Item(19, 5).Value = 1541271038.75871
$ws.Cells.Item(19, 6).Value = 0.004131978223816866
$ws.Cells.Item(19, 7).Value = 26.49260867365715
$ws.Cells.Item(19, 8).Value = 20
$ws.Cells.Item(19, 9).Value = $true
$ws.Cells.Item(19, 10).Value = 76923.07692307692
$ws.Cells.Item(19, 11).Value = 108270.0903728412
$ws.Cells.Item(19, 12).Value = 18.32631702039958
$ws.Cells.Item(19, 13).Value = 44.81892569405673
$ws.Cells.Item(20, 3).Value = 165.4740402215803
$ws.Cells.Item(20, 4).Value = 80000888
$ws.Cells.Item(20, 5).Value = 1386432394.135247
$ws.Cells.Item(20, 6).Value = 0.002262925582568383
$ws.Cells.Item(20, 7).Value = 47.74149181262842
$ws.Cells.Item(20, 8).Value = 15
$ws.Cells.Item(20, 9).Value = $true
$ws.Cells.Item(20, 10).Value = 76923.07692307692
$ws.Cells.Item(20, 11).Value = 71532.30099033478
$ws.Cells.Item(20, 12).Value = 27.73840590236427
$ws.Cells.Item(20, 13).Value = 75.47989771499269
$ws.Cells.Item(21, 3).Value = 71.22950717745928
$ws.Cells.Item(21, 4).Value = 80230776
$ws.Cells.Item(21, 5).Value = 427332063.0859565
$ws.Cells.Item(21, 6).Value = 0.003963270018734173
$ws.Cells.Item(21, 7).Value = 66.86601742069158
$ws.Cells.Item(21, 8).Value = 9
$ws.Cells.Item(21, 9).Value = $false
$ws.Cells.Item(21, 10).Value = 76923.07692307692
$ws.Cells.Item(21, 11).Value = 105410.4746201164
$ws.Cells.Item(21, 12).Value = 18.82348037185803
$ws.Cells.Item(21, 13).Value = 85.6894977925496
